$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 08:20:27"
$wsOverview.Range("G5").Value = "2016-09-01 08:20:27"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-01 08:20:16"
$wsZhCn.Range("H5").Value = "2016-09-01 08:20:16"
$wsZhCn.Range("K2").Value = "2016-09-01 08:20:52"
$wsZhCn.Range("K5").Value = "2016-09-01 08:20:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-01 08:20:27"
$wsDeDe.Range("H5").Value = "2016-09-01 08:20:27"
$wsDeDe.Range("K2").Value = "2016-09-01 08:20:59"
$wsDeDe.Range("K5").Value = "2016-09-01 08:20:59"
